$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "csAUq754"
$ws.Range("B2").Value = 23072557
$ws.Range("C2").Value = "jtyhqyk59"
$ws.Range("D2").Value = "r38!vMV$"
$ws.Range("F2").Value = "QJgaITVt"
$ws.Range("G2").Value = "tLjA"

$wb.Save()
